$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 23:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 362573
$ws.Range("C4").Value = 25900
$ws.Range("E4").Value = 332540
$ws.Range("G4").Value = 1104
$ws.Range("H4").Value = 10720

# Row 5 - Espana
$ws.Range("B5").Value = 136675
$ws.Range("C5").Value = 5029
$ws.Range("E5").Value = 82897
$ws.Range("G5").Value = 700
$ws.Range("H5").Value = 13341

# Row 7 - Francia
$ws.Range("B7").Value = 102453
$ws.Range("C7").Value = 2330
$ws.Range("E7").Value = 72018
$ws.Range("F7").Value = 4895
$ws.Range("G7").Value = 151
$ws.Range("H7").Value = 1735

# Row 16 - Canada
$ws.Range("B16").Value = 16558
$ws.Range("C16").Value = 1046
$ws.Range("D16").Value = 3534
$ws.Range("E16").Value = 12702
$ws.Range("G16").Value = 42
$ws.Range("H16").Value = 322

# Row 153 - Eritrea
$ws.Range("B153").Value = 31
$ws.Range("C153").Value = 2
$ws.Range("E153").Value = 31
